# Actualización de documentos, Semana 4
#
# The "Ultima Revisión" status table (shape "Tabla 27") bumps its two
# week counters by one:
#   Documentos : 13 -> 14
#   Aplicación : 14 -> 15
#
# Also refreshes the datetimeFigureOut footer field (9/4/2013 -> 9/5/2013)
# wherever it appears across the deck (most slides carry it in their
# footer placeholder; slides without one are simply skipped).

$p = $ppt.ActivePresentation

function Update-RevisionTable($tbl) {
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
            $cellShape = $tbl.Cell($r, $c).Shape
            $txt = $cellShape.TextFrame.TextRange.Text
            if ($txt -eq "13") {
                $cellShape.TextFrame.TextRange.Text = "14"
            } elseif ($txt -eq "14") {
                $cellShape.TextFrame.TextRange.Text = "15"
            }
        }
    }
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)

        if ($shp.HasTable -and $shp.Name -eq "Tabla 27") {
            Update-RevisionTable $shp.Table
        }

        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "9/4/2013") {
                $tr.Text = "9/5/2013"
            }
        }
    }
}
